# Insert two new weekly price rows for "Espinaca" (Mercado Mayorista Lo Valledor
# de Santiago) right before the existing row 359, pushing all subsequent rows
# down by two. This reproduces the upstream diff, which shows every row from
# 359 through the previous last row (421) shifting down to 361..423, with the
# content for the two brand-new rows (new 359 and new 360) being the data
# below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 359 (shifts old row 359 onward down to 361 onward)
$ws.Range("A359:A360").EntireRow.Insert()

# --- New row 359 ---
$ws.Range("A359").Value2 = 6
$ws.Range("B359").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C359").Value2 = "Metropolitana"
$ws.Range("D359").Value2 = 44522
$ws.Range("E359").Value2 = 13
$ws.Range("F359").Value2 = 100112012
$ws.Range("G359").Value2 = "Espinaca"
$ws.Range("H359").Value2 = "Sin especificar"
$ws.Range("I359").Value2 = "Primera"
$ws.Range("J359").Value2 = 210
$ws.Range("K359").Value2 = 5500
$ws.Range("L359").Value2 = 6000
$ws.Range("M359").Value2 = 5786
$ws.Range("N359").Value2 = "`$/cuna 10 kilos"
$ws.Range("O359").Value2 = "Provincia de Chacabuco"
$ws.Range("P359").Value2 = 579
$ws.Range("Q359").Value2 = 10
$ws.Range("R359").Value2 = "Hortaliza"

# --- New row 360 ---
$ws.Range("A360").Value2 = 6
$ws.Range("B360").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C360").Value2 = "Metropolitana"
$ws.Range("D360").Value2 = 44522
$ws.Range("E360").Value2 = 13
$ws.Range("F360").Value2 = 100112012
$ws.Range("G360").Value2 = "Espinaca"
$ws.Range("H360").Value2 = "Sin especificar"
$ws.Range("I360").Value2 = "Primera"
$ws.Range("J360").Value2 = 280
$ws.Range("K360").Value2 = 5500
$ws.Range("L360").Value2 = 6000
$ws.Range("M360").Value2 = 5732
$ws.Range("N360").Value2 = "`$/cuna 10 kilos"
$ws.Range("O360").Value2 = "Región Metropolitana"
$ws.Range("P360").Value2 = 573
$ws.Range("Q360").Value2 = 10
$ws.Range("R360").Value2 = "Hortaliza"

Write-Host "Inserted rows 359-360; new dimension:" $ws.UsedRange.Address()
